# QA Compiler: Fix DAILY delta calculation for John's row (01/06) and TOTAL row,
# plus corresponding chart-data tables lower on the DAILY sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAILY")

# Row 7 (01/06) - John had no actual delta today, so his stats become "--"
$ws.Range("J7").Value = "--"
$ws.Range("K7").Value = "--"
# Force text storage for the percent-looking label so Excel doesn't
# auto-convert it into a numeric percentage value.
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "0.0%"
$ws.Range("M7").Value = "--"
$ws.Range("Q7").Value = "--"

# Row 8 (TOTAL) - recompute totals now that John's daily delta is 0
$ws.Range("J8").Value = 4
$ws.Range("K8").Value = 2
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "40.0%"
$ws.Range("Q8").Value = 5

# Chart data table (Comp % chart) - John's 01/06 value drops to 0
$ws.Range("D14").Value = 0

# Chart data table (Issues/Pending chart) - John's 01/06 value drops to 0
$ws.Range("D35").Value = 0
